$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / tab to reflect the new "through" date
$ws.Name = "Through 2022-07-06"

# Update the row label for July to reflect the new "through" date
$ws.Range("A8").Value = "July (through 07-06)"

# Update the July row (row 8) figures for the years that changed
$ws.Cells.Item(8, 3).Value = 9    # C8 (2016)
$ws.Cells.Item(8, 5).Value = 17   # E8 (2018)
$ws.Cells.Item(8, 7).Value = 20   # G8 (2020)
$ws.Cells.Item(8, 8).Value = 29   # H8 (2021)
$ws.Cells.Item(8, 9).Value = 30   # I8 (2022)

# Update the Total row (row 9) figures for the years that changed
$ws.Cells.Item(9, 3).Value = 257  # C9 (2016)
$ws.Cells.Item(9, 5).Value = 370  # E9 (2018)
$ws.Cells.Item(9, 7).Value = 492  # G9 (2020)
$ws.Cells.Item(9, 8).Value = 789  # H9 (2021)
$ws.Cells.Item(9, 9).Value = 836  # I9 (2022)
